# API: Gameweeks import (#25)
# Adds two new columns ("Show Statistics Continuously", "Gameweek") to the
# "Challenges" sheet, with sample data "true" / 1 in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New header cells
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# New sample-data cells in row 2.
# S2 must stay a *text* string "true" (not a boolean) - the leading
# apostrophe forces Excel to treat it as text rather than auto-converting
# to a boolean. Resetting the style afterwards drops the "quote prefix"
# formatting flag that the apostrophe trick leaves behind, so the cell
# keeps a plain/default style (matching a plain text import).
$sCell = $ws.Range("S2")
$sCell.Value = "'true"
$sCell.Style = "Normal"

$ws.Range("T2").Value = 1
